$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.920.79"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "1.819.53"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'309.65"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4677"
$ws.Range("E7").Value = "  +2.99%  "
$ws.Range("D8").Value = "'0.3697"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "'0.07383"
$ws.Range("E9").Value = "  +2.52%  "
$ws.Range("D10").Value = "'0.8717"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("D11").Value = "'20.47"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "1.847.87"
$ws.Range("E12").Value = "  +3.38%  "
$ws.Range("D13").Value = "'5.366"
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("D14").Value = "'92.44"
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").Value = "'0.07073"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "'6.499"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "'0.000008724"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'14.78"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "26.957.72"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "'5.351"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").Value = "'10.60"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").Value = "2.043.33"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").Value = "'1.904"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "'151.41"
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("D27").Value = "'2.182"
$ws.Range("E27").Value = "  +2.45%  "
$ws.Range("E28").Value = "  +2.39%  "
$ws.Range("D29").Value = "'5.332"
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").Value = "'115.80"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("D31").Value = "'0.08942"
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("D32").Value = "'0.7699"
$ws.Range("E32").Value = "  +2.73%  "
$ws.Range("D33").Value = "'1.166"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").Value = "'4.512"
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("D35").Value = "'2.901"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("D38").Value = "'0.01970"
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("D39").Value = "'0.05291"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").Value = "'7.314"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("D41").Value = "'2.949"
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D42").Value = "'0.5345"
$ws.Range("D43").Value = "'2.361"
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("D44").Value = "'0.1669"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("D45").Value = "'8.458"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "'0.4961"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("D47").Value = "'10.44"
$ws.Range("E47").Value = "  +2.22%  "
$ws.Range("D48").Value = "'104.23"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.674"
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "'0.06283"
